$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-354).
# All of them are being bumped forward by one day: 46081 -> 46082.
$range = $ws.Range("C2:C354")
$range.Value = 46082
